# Trade #32 closed at 2026-02-16 21:28:22 - leadlag DOWN +0.000%
#
# This script applies the following changes to the workbook:
#  1. Summary sheet: refresh OVERALL and leadlag aggregate stats.
#  2. leadlag sheet: close out trade #14 (row 13) that was previously OPEN,
#     and append the newly-opened trade #32 as a new row (28).
#  3. All Trades sheet: append the same closed trade as a new row (15).
#  4. Comparison sheet: refresh leadlag aggregate stats.

$wb = $excel.ActiveWorkbook

# Excel's COM layer auto-detects numbers/dates/percentages typed into a
# cell and silently converts them (e.g. "64.3%" -> 0.643, "2026-02-16" ->
# a date serial). The source values must stay literal text, so for those
# ambiguous strings we force the cell to Text format first.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 14
Set-TextValue $summary.Range("D2") "64.3%"
Set-TextValue $summary.Range("E2") "+2.6297%"
Set-TextValue $summary.Range("F2") "+0.1878%"

$summary.Range("C3").Value = 26
Set-TextValue $summary.Range("D3") "30.8%"
Set-TextValue $summary.Range("E3") "+2.5844%"
Set-TextValue $summary.Range("F3") "+0.0994%"

# ---------------------------------------------------------------------------
# 2. leadlag sheet
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

# Close out existing trade in row 13 (Trade # 14)
$leadlag.Range("G13").Value = 69115.625505
$leadlag.Range("H13").Value = "CLOSED"
$leadlag.Range("I13").Value = 0.3847
$leadlag.Range("J13").Value = 3.85
$leadlag.Range("M13").Value = "time_exit_5min"
$leadlag.Range("N13").Value = 5

# Append new trade (Trade # 32) as row 28
$leadlag.Range("A28").Value = 32
Set-TextValue $leadlag.Range("B28") "2026-02-16"
$leadlag.Range("C28").Value = "21:28:22"
$leadlag.Range("D28").Value = "leadlag"
$leadlag.Range("E28").Value = "DOWN"
$leadlag.Range("F28").Value = 68889.565
$leadlag.Range("H28").Value = "OPEN"
$leadlag.Range("I28").Value = 0
$leadlag.Range("J28").Value = 0
$leadlag.Range("K28").Value = 0.7052
$leadlag.Range("L28").Value = "Coinbase leading with -0.071% move"
$leadlag.Range("N28").Value = 0

# ---------------------------------------------------------------------------
# 3. All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A15").Value = 14
Set-TextValue $allTrades.Range("B15") "2026-02-16"
$allTrades.Range("C15").Value = "21:23:08"
$allTrades.Range("D15").Value = "leadlag"
$allTrades.Range("E15").Value = "DOWN"
$allTrades.Range("F15").Value = 69382.565
$allTrades.Range("G15").Value = 69115.625505
$allTrades.Range("H15").Value = "CLOSED"
$allTrades.Range("I15").Value = 0.3847
$allTrades.Range("J15").Value = 3.85
$allTrades.Range("K15").Value = 0.6303
$allTrades.Range("L15").Value = "Binance leading with -0.063% move"
$allTrades.Range("M15").Value = "time_exit_5min"
$allTrades.Range("N15").Value = 5

# ---------------------------------------------------------------------------
# 4. Comparison sheet
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 26
Set-TextValue $comparison.Range("C2") "30.8%"
Set-TextValue $comparison.Range("D2") "3.22"
Set-TextValue $comparison.Range("E2") "+0.4688%"
Set-TextValue $comparison.Range("G2") "1.61"
